$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45243
}
